$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row to match the new "import exercise data" API field
# names (Description/Instructions/Image/Video/Muscle/Equipment columns ->
# Category/Instructions/imageUrl/videoUrl/muscleGroups/Equipment).
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Instructions"
$ws.Range("D1").Value = "imageUrl"
$ws.Range("E1").Value = "videoUrl"
$ws.Range("F1").Value = "muscleGroups"
$ws.Range("G1").Value = "Equipment"

# Move the active selection (matches the saved worksheet view after editing).
$ws.Range("B31").Select()

# Ensure the sheet prints in portrait orientation.
$ws.PageSetup.Orientation = 1
